# Automatische test-sync: 2025-08-05 18:48:50
#
# Adds the newest "Testmail #20" log entry to the Logs sheet, extends the
# conditional-formatting ranges to cover the new row, and refreshes the
# category counts on the Dashboard sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Logs sheet: append the new row (row 41)
# ---------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$newRow = 41

$logs.Cells.Item($newRow, 1).Value  = "Ik ben niet tevreden over hoe dit is gegaan."
$logs.Cells.Item($newRow, 2).Value  = "mailmind.test@zohomail.eu"
$logs.Cells.Item($newRow, 3).Value  = "Testmail #20: Ik ben niet tevreden over hoe dit is gegaan."
$logs.Cells.Item($newRow, 4).Value  = "Klacht / Probleem"
$logs.Cells.Item($newRow, 5).Value  = "Bedankt, we hebben dit doorgestuurd naar klachten@bedrijf.nl."
$logs.Cells.Item($newRow, 6).Value  = "2025-08-05 18:48:05"
$logs.Cells.Item($newRow, 7).Value  = "Ja"
$logs.Cells.Item($newRow, 8).Value  = "Ja"
$logs.Cells.Item($newRow, 9).Value  = "Nee"
$logs.Cells.Item($newRow, 10).Value = "Nee"

# ---------------------------------------------------------------------
# 2) Logs sheet: extend the conditional-formatting ranges from row 40
#    to row 41 for every formatted column.
# ---------------------------------------------------------------------
$ccols = @("D", "G", "H", "I", "J")
foreach ($col in $ccols) {
    $oldRange = $logs.Range("$col" + "2:" + "$col" + "40")
    $newRange = $logs.Range("$col" + "2:" + "$col" + "41")
    for ($i = 1; $i -le $oldRange.FormatConditions.Count; $i++) {
        $oldRange.FormatConditions.Item($i).ModifyAppliesToRange($newRange)
    }
}

# ---------------------------------------------------------------------
# 3) Dashboard sheet: refresh the category counts now that the new
#    "Klacht / Probleem" entry moved it ahead of "Opvolging / Status".
# ---------------------------------------------------------------------
$dashboard = $wb.Worksheets.Item("Dashboard")

$dashboard.Cells.Item(6, 1).Value = "Klacht / Probleem"
$dashboard.Cells.Item(6, 2).Value = 2

$dashboard.Cells.Item(7, 1).Value = "Opvolging / Status"
$dashboard.Cells.Item(7, 2).Value = 2
